$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 19:25"

# --- Country name re-mapping for rows 146-148 (Botsuana moved up in source list) ---
# Row 146 was Guyana -> now Botsuana (with fresh Botsuana stats)
# Row 147 was Mali   -> now Guyana   (keeps former Guyana stats)
# Row 148 was Botsuana -> now Mali   (keeps former Mali stats)
# Row 149 stays Sudan del Sur (stats refreshed)
$ws.Range("A146").Value = "Botsuana"
$ws.Range("A147").Value = "Guyana"
$ws.Range("A148").Value = "Mali"

# --- Updated statistics (Casos totales / Nuevos casos / Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes) ---
# Row 4
$ws.Range("B4").Value = 8003660
$ws.Range("C4").Value = 11662
$ws.Range("D4").Value = 5138536
$ws.Range("E4").Value = 2645319
$ws.Range("G4").Value = 110
$ws.Range("H4").Value = 219805

# Row 5
$ws.Range("B5").Value = 7172559
$ws.Range("C5").Value = 53259
$ws.Range("D5").Value = 6223372
$ws.Range("E5").Value = 839306
$ws.Range("G5").Value = 697
$ws.Range("H5").Value = 109881

# Row 25
$ws.Range("B25").Value = 329273
$ws.Range("C25").Value = 2982
$ws.Range("E25").Value = 42662
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 9711

# Row 27
$ws.Range("B27").Value = 293553
$ws.Range("C27").Value = 3060
$ws.Range("D27").Value = 238643
$ws.Range("E27").Value = 52894
$ws.Range("G27").Value = 36
$ws.Range("H27").Value = 2016

# Row 33
$ws.Range("B33").Value = 153761
$ws.Range("C33").Value = 1357
$ws.Range("D33").Value = 129498
$ws.Range("E33").Value = 21627
$ws.Range("G33").Value = 31
$ws.Range("H33").Value = 2636

# Row 62
$ws.Range("E62").Value = 124
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 28

# Row 65
$ws.Range("B65").Value = 54624
$ws.Range("C65").Value = 1056
$ws.Range("D65").Value = 23941
$ws.Range("E65").Value = 30217
$ws.Range("G65").Value = 7
$ws.Range("H65").Value = 466

# Row 66
$ws.Range("B66").Value = 53325
$ws.Range("C66").Value = 253
$ws.Range("D66").Value = 37382
$ws.Range("E66").Value = 14134
$ws.Range("G66").Value = 8
$ws.Range("H66").Value = 1809

# Row 108
$ws.Range("B108").Value = 10180
$ws.Range("C108").Value = 10
$ws.Range("D108").Value = 9855
$ws.Range("E108").Value = 256

# Row 135
$ws.Range("B135").Value = 4844
$ws.Range("C135").Value = 92
$ws.Range("E135").Value = 1514

# Row 146
$ws.Range("B146").Value = 3515
$ws.Range("C146").Value = 296
$ws.Range("D146").Value = 853
$ws.Range("E146").Value = 2642
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = 20

# Row 147
$ws.Range("B147").Value = 3469
$ws.Range("D147").Value = 2318
$ws.Range("E147").Value = 1048
$ws.Range("H147").Value = 103

# Row 148
$ws.Range("B148").Value = 3286
$ws.Range("D148").Value = 2527
$ws.Range("E148").Value = 627
$ws.Range("H148").Value = 132

# Row 149
$ws.Range("B149").Value = 2787
$ws.Range("C149").Value = 10
$ws.Range("E149").Value = 1442

